# Refresh the crypto price/volume table per the commit diff.
#
# Columns B-E are all stored as text (inline strings) in this workbook.
# A plain Range.Value assignment auto-converts numeric-looking strings
# (e.g. "1.00", "0.574", "7.08") into real numbers, which would silently
# drop the literal formatting (trailing zeros, etc.). Set-TextValue below
# detects that case and forces the cell to Text format first so the exact
# string survives, then restores the default "Normal" style so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($RangeAddr, $Text)
    $cell = $ws.Range($RangeAddr)
    if ($Text -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $cell.NumberFormat = "@"
        $cell.Value = $Text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $Text
    }
}

Set-TextValue "D2" "65.551.00"
Set-TextValue "E2" "  +0.04%  "
Set-TextValue "D3" "3.561.14"
Set-TextValue "E3" "  +3.58%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "604.91"
Set-TextValue "E5" "  +2.06%  "
Set-TextValue "D6" "141.32"
Set-TextValue "E6" "  +4.07%  "
Set-TextValue "D7" "3.559.12"
Set-TextValue "E7" "  +3.55%  "
Set-TextValue "E8" "  +0.09%  "
Set-TextValue "E9" "  +0.27%  "
Set-TextValue "E10" "  +3.14%  "
Set-TextValue "D11" "7.08"
Set-TextValue "E11" "  -3.81%  "
Set-TextValue "E12" "  +4.89%  "
Set-TextValue "D13" "4.166.68"
Set-TextValue "E13" "  +3.72%  "
Set-TextValue "E14" "  +4.54%  "
Set-TextValue "D15" "27.33"
Set-TextValue "E15" "  +2.64%  "
Set-TextValue "D16" "3.557.53"
Set-TextValue "E16" "  +3.28%  "
Set-TextValue "E17" "  +1.63%  "
Set-TextValue "D18" "65.558.55"
Set-TextValue "E18" "  +0.14%  "
Set-TextValue "E19" "  +4.75%  "
Set-TextValue "E20" "  +1.76%  "
Set-TextValue "D21" "14.36"
Set-TextValue "E21" "  +5.28%  "
Set-TextValue "D22" "396.06"
Set-TextValue "E22" "  +0.39%  "
Set-TextValue "D23" "0.574"
Set-TextValue "E23" "  +5.02%  "
Set-TextValue "D24" "3.702.48"
Set-TextValue "E24" "  +3.45%  "
Set-TextValue "D25" "74.23"
Set-TextValue "E25" "  +1.22%  "
Set-TextValue "E26" "  +0.01%  "
Set-TextValue "D27" "0.0000118"
Set-TextValue "E27" "  +11.76%  "
Set-TextValue "D28" "7.91"
Set-TextValue "E28" "  +9.16%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  -0.14%  "
Set-TextValue "E30" "  +2.11%  "
Set-TextValue "E31" "  +3.04%  "
Set-TextValue "D32" "3.574.37"
Set-TextValue "E32" "  +3.83%  "
Set-TextValue "B33" "Kaspa"
Set-TextValue "C33" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D33" "0.148"
Set-TextValue "E33" "  +0.80%  "
Set-TextValue "B34" "USDe"
Set-TextValue "C34" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D34" "1.00"
Set-TextValue "E34" "  +0.06%  "
Set-TextValue "D35" "23.90"
Set-TextValue "E35" "  +3.86%  "
Set-TextValue "E36" "  +6.90%  "
Set-TextValue "D37" "7.09"
Set-TextValue "E37" "  +2.30%  "
Set-TextValue "E38" "  +3.60%  "
Set-TextValue "D39" "168.21"
Set-TextValue "E39" "  -1.72%  "
Set-TextValue "D40" "5.06"
Set-TextValue "E40" "  +5.01%  "
Set-TextValue "D41" "0.0814"
Set-TextValue "E41" "  +5.61%  "
Set-TextValue "D42" "0.836"
Set-TextValue "E42" "  +1.59%  "
Set-TextValue "D43" "26.49"
Set-TextValue "E43" "  +16.93%  "
Set-TextValue "D44" "43.16"
Set-TextValue "E44" "  -0.80%  "
Set-TextValue "E45" "  +0.02%  "
Set-TextValue "D46" "4.47"
Set-TextValue "E46" "  +1.00%  "
Set-TextValue "B47" "ONDO"
Set-TextValue "C47" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D47" "1.22"
Set-TextValue "E47" "  +10.37%  "
Set-TextValue "B48" "Stacks"
Set-TextValue "C48" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D48" "1.71"
Set-TextValue "E48" "  +5.92%  "
Set-TextValue "D49" "2.486.64"
Set-TextValue "E49" "  +13.01%  "
Set-TextValue "D50" "6.85"
Set-TextValue "E50" "  +4.50%  "
Set-TextValue "D51" "2.37"
Set-TextValue "E51" "  +19.11%  "
